$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The contract "009/PV009" (row 2) was resiliated/cancelled and must be
# removed from the situation. Delete the whole row; subsequent rows shift
# up by one.
$ws.Rows(2).Delete()

# Update the totals row (now row 8) to account for the removed contract's
# contribution (avance=10000, caution=10000, net=10000).
$ws.Range("I8").Value = 17500
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 33550
